# Added scrolling text output for testing
# Populate the "Wildcard Ranges" column (I) for the existing Typography rows
# (4-9) with the byte range used by the new scrolling-text test entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Typography")

$ws.Range("I4").Value = "0x00-0xFF"
$ws.Range("I5").Value = "0x00-0xFF"
$ws.Range("I6").Value = "0x00-0xFF"
$ws.Range("I7").Value = "0x00-0xFF"
$ws.Range("I8").Value = "0x00-0xFF"
$ws.Range("I9").Value = "0x00-0xFF"
